# Time Tracking update 5/15
# updated last weeks and what I worked on Sunday/Today

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sunday (row 9, "Week 6") hours worked: 3.5 -> 10.5
$ws.Range("H9").Value = 10.5

# Today (row 10, "Week 7") hours worked: was blank, now 4
$ws.Range("H10").Value = 4

# Move the active selection from H10 to H11 to reflect where the user left off
[void]$ws.Range("H11").Select()
